$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C entirely (shrinks used range from A1:C16 to A1:B16)
$ws.Range("C1:C16").Delete() | Out-Null

# Update column B values with the new simulation results
$ws.Range("B1").Value = 1

$ws.Range("B2").Value = 0.330175302796987
$ws.Range("B3").Value = 37.67534303548216
$ws.Range("B4").Value = 942.5614073070984
$ws.Range("B5").Value = 108.4618374190834
$ws.Range("B6").Value = 52402.88008616195
$ws.Range("B7").Value = 2454.100102838397
$ws.Range("B8").Value = 1992.735548603491
$ws.Range("B9").Value = 418.7325257263029
$ws.Range("B10").Value = 2026.407514986833
$ws.Range("B11").Value = 2941.891519169702
$ws.Range("B12").Value = -72.29922816759412
$ws.Range("B13").Value = -1.664808906372224
$ws.Range("B14").Value = -1.611548195128604
$ws.Range("B15").Value = 1.772852097089415
$ws.Range("B16").Value = 2
